$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated CE (column F) values - diameter now measured by max side length of cell
$updates = @{
    20 = -3.036088916670467
    21 = -3.049772684337528
    22 = -3.062305129253268
    23 = -3.074837574169008
    24 = -3.086218696333427
    25 = -3.097599818497844
    27 = -3.118059417324032
    28 = -3.128289216737126
    32 = -3.150736926868178
    33 = -3.15634885440094
    35 = -3.164752556196468
    36 = -3.17177456036423
    37 = -3.174566334626995
    39 = -3.181559959787522
    40 = -3.173178683114247
    41 = -3.191373738218049
    42 = -3.195575589115813
    43 = -3.182178210893511
    44 = -3.183159685670903
    45 = -3.208181141809104
    46 = -3.192749481960834
    47 = -3.200395691119567
    48 = -3.184964031271296
    50 = -3.170914069050793
    51 = -3.19333935405959
    52 = -3.165324566640283
    53 = -3.16670626426805
    54 = -3.165267808625817
    55 = -3.168059582888583
    56 = -3.154037999675312
    57 = -3.138606339827043
    81 = -3.314660256527335
    82 = -3.35733023067796
    83 = -3.374424450610687
    84 = -3.367895768396122
    85 = -3.33681118517157
    86 = -3.292872699282523
    87 = -3.238807390983804
    88 = -3.166436068057369
    89 = -3.07771351810183
    92 = -3.403462648016482
    93 = -3.428779909822958
    94 = -3.453235465287177
    95 = -3.448114041105835
    97 = -3.380209903762774
    99 = -3.274345801950522
    103 = -3.484044852034426
    105 = -3.511793772955028
    106 = -3.510859999399067
    107 = -3.482863411205783
    108 = -3.443640925868633
    109 = -3.390580378189961
    111 = -3.253807200636091
    114 = -3.523966434901569
    116 = -3.55029254616861
    117 = -3.541638145036337
    119 = -3.480636100333307
    120 = -3.431651437861417
    122 = -3.304415182263253
    126 = -3.582113035184955
    127 = -3.584759037260379
    128 = -3.576956417286147
    129 = -3.550782874948943
    130 = -3.51683158513027
    131 = -3.46775189936583
    132 = -3.406829592405015
    133 = -3.334430086819109
    136 = -3.591117066158365
    138 = -3.604199530225025
    139 = -3.591832354253244
    140 = -3.566919274564132
    141 = -3.534303096331451
    142 = -3.488656776629619
    144 = -3.361693720787636
    147 = -3.614620966770244
    149 = -3.626436992406192
    150 = -3.611445529703895
    152 = -3.555050621115603
    153 = -3.509118804928575
    154 = -3.452111825512828
    158 = -3.632305967427921
    159 = -3.639932681903296
    160 = -3.638739201579023
    162 = -3.602688293622485
    163 = -3.570561774151418
    164 = -3.523288573677846
    166 = -3.40014629304635
    171 = -3.653846701045347
    173 = -3.618224206723761
    174 = -3.582565485642089
    176 = -3.481473337760865
    177 = -3.414439095769947
    180 = -3.660211462224339
    182 = -3.664177105565324
    183 = -3.650791699369308
    186 = -3.546589553977225
    187 = -3.491913827901804
    188 = -3.425896289784562
    191 = -3.672907202617076
    192 = -3.679459286447214
    193 = -3.674319265515693
    194 = -3.661415195237005
    195 = -3.637584380543003
    196 = -3.60213297547393
    197 = -3.556561381711564
    198 = -3.501104921957951
    199 = -3.436297411727539
    202 = -3.682128376707835
    203 = -3.687535923126708
    204 = -3.683780518249396
    205 = -3.667949197147848
    206 = -3.643349315096021
    209 = -3.508693314253891
    210 = -3.445331569455059
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

Write-Output "Updated $($updates.Count) cells in column F"